# Fruta / hortaliza, semanal
# Insert one new weekly record before the current row 343, shifting all
# subsequent rows (343-367) down by one position (to 344-368), then
# populate the newly inserted row with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 343; this pushes the former rows
# 343..367 down to 344..368, preserving all of their data untouched.
$ws.Rows("343").Insert()

# Populate the newly inserted row 343 with the new weekly record.
$ws.Cells.Item(343, 1).Value = 7
$ws.Cells.Item(343, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(343, 3).Value = "Ñuble"
$ws.Cells.Item(343, 4).Value = 45223
$ws.Cells.Item(343, 5).Value = 16
$ws.Cells.Item(343, 6).Value = 100112032
$ws.Cells.Item(343, 7).Value = "Zapallo italiano"
$ws.Cells.Item(343, 8).Value = "Sin especificar"
$ws.Cells.Item(343, 9).Value = "Primera"
$ws.Cells.Item(343, 10).Value = 100
$ws.Cells.Item(343, 11).Value = 13000
$ws.Cells.Item(343, 12).Value = 14000
$ws.Cells.Item(343, 13).Value = 13500
$ws.Cells.Item(343, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(343, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(343, 16).Value = 270
$ws.Cells.Item(343, 17).Value = 50
$ws.Cells.Item(343, 18).Value = "Hortaliza"
